# Slide 17 ("Fazit – Potentiale"), shape "Textplatzhalter 11":
#   - append a period to the end of the "Hinzufügen ..." bullet
#   - append a period to the end of the "Eingabe ..." bullet
#   - reword the "Ausgabe ..." bullet and give it a trailing period
#
# We edit via substring Characters() ranges (rather than re-assigning the
# whole TextRange.Text) so that each run's existing character formatting
# (rPr: lang/sz/etc.) is preserved and only the text content changes.

function Replace-Substring {
    param($TextRange, [string]$OldText, [string]$NewText)

    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Could not find expected text: '$OldText'"
    }
    $sub = $TextRange.Characters($idx + 1, $OldText.Length)
    $sub.Text = $NewText
}

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(17)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

Replace-Substring $tr "Hinzufügen von neuen Kampagnen-Ergebnissen" "Hinzufügen von neuen Kampagnen-Ergebnissen."
Replace-Substring $tr "Eingabe von angebotenen Liefermengen" "Eingabe von angebotenen Liefermengen."
Replace-Substring $tr "Ausgabe von Budget-Allokation und Umsatz-Prognose" "Ausgabe von Budget-Allokation für den maximal erzielbaren Umsatz."
